$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'326.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.21%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.30%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.528"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-4.36%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08072"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.98%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.705"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.23%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.340"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.54%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.898"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.19%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-7.30%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9477"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.48%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1181"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.88%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1897"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.51%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.22%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.04176"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'5.55%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1065"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.08%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001271"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.65%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005992"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.19%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.601"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.59%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.65%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.420"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-7.31%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1374"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.18%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2533"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.53%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04253"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-3.49%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-1.22%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004504"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.79%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001233"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.57%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004000"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.21%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02654"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-6.30%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05537"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.76%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007709"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.36%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1393"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.32%"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "Dexo"
$ws.Range("C43").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D43").Value = "'0.006670"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-26.52%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002060"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.99%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008684"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-16.61%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007116"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.27%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.20%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003423"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-13.75%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002276"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.20%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.20%"
$ws.Range("E51").Style = "Normal"
